$wb = $excel.ActiveWorkbook

# ----- Sheet "Horario" (sheet1) -----
$ws1 = $wb.Worksheets.Item("Horario")

# 9:00-10:00 (rows 2-3), column B: keep only "Introducción a la Programación (1) / Reloj 102"
$ws1.Range("B2:B3").Value = "Introducción a la Programación (1)`nReloj 102"

# 11:00-12:00 (rows 4-5), column B: now holds "Introducción a la Programación (2) / Reloj 102"
$ws1.Range("B4:B5").Value = "Introducción a la Programación (2)`nReloj 102"

# 14:00-16:00 (rows 7-9): remove the old "Algoritmos y Complejidad (3) / Reloj 103" from column B
$ws1.Range("B7:B9").ClearContents()

# 14:00-16:00 (rows 7-9), column F: "Algoritmos y Complejidad (3) / Ciencias 507"
# (this replaces the old "Introducción a la Programación (asdf) / Ciencias 507" that used to
# occupy F8:F9 and extends the entry to F7 as well)
$ws1.Range("F7:F9").Value = "Algoritmos y Complejidad (3)`nCiencias 507"

# ----- Sheet "Tabla" (sheet2) -----
$ws2 = $wb.Worksheets.Item("Tabla")

# Sección 2 now runs 11:00-13:00 in Reloj 102 (was 9:00-11:00 in Reloj 103)
$ws2.Range("C3").Value = "11:00"
$ws2.Range("D3").Value = "13:00"
$ws2.Range("E3").Value = "Reloj 102"

# Sección 3 now happens on Viernes in Ciencias 507 (was Lunes in Reloj 103)
$ws2.Range("B4").Value = "Viernes"
$ws2.Range("E4").Value = "Ciencias 507"

# The "Sección asdf" row (row 5) is removed entirely
$ws2.Rows.Item(5).Delete() | Out-Null
